# Update countries & provincias Spain
# - Swap rank order of Niger / Burkina Faso (rows 94/95) with refreshed stats
# - Swap rank order of Martinica / Guatemala (rows 122/123) with refreshed stats
# - Update several countries' numeric stats (rows 4, 16, 21, 57, 81, 179)
# - Update the "Datos actualizados" timestamp string

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp update (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 12 de Abril de 2020 a las 20:22"

# --- Row 4 (Estados Unidos) ---
$ws.Range("B4").Value = 550655
$ws.Range("C4").Value = 17776
$ws.Range("E4").Value = 497868
$ws.Range("F4").Value = 11760
$ws.Range("G4").Value = 1090
$ws.Range("H4").Value = 21667

# --- Row 16 (Canada) ---
$ws.Range("B16").Value = 24283
$ws.Range("C16").Value = 965
$ws.Range("D16").Value = 7106
$ws.Range("E16").Value = 16464
$ws.Range("G16").Value = 60
$ws.Range("H16").Value = 713

# --- Row 21 (Israel) ---
$ws.Range("B21").Value = 11145
$ws.Range("C21").Value = 402
$ws.Range("D21").Value = 1627
$ws.Range("E21").Value = 9415
$ws.Range("F21").Value = 183

# --- Row 57 (Egipto) ---
$ws.Range("B57").Value = 2065
$ws.Range("C57").Value = 126
$ws.Range("D57").Value = 589
$ws.Range("E57").Value = 1317
$ws.Range("G57").Value = 13
$ws.Range("H57").Value = 159

# --- Row 81 (Tunez) ---
$ws.Range("B81").Value = 707
$ws.Range("C81").Value = 22
$ws.Range("E81").Value = 633
$ws.Range("G81").Value = 3
$ws.Range("H81").Value = 31

# --- Rows 94/95: Niger / Burkina Faso swap ranking ---
# Burkina Faso overtakes Niger: new row 94 = Burkina Faso (fresh numbers),
# new row 95 = Niger (carries forward the old row-94 numbers unchanged).
$ws.Range("A94").Value = "Burkina Faso"
$ws.Range("B94").Value = 497
$ws.Range("C94").Value = 13
$ws.Range("D94").Value = 161
$ws.Range("E94").Value = 309
$ws.Range("H94").Value = 27

$ws.Range("A95").Value = "Niger"
$ws.Range("B95").Value = 491
$ws.Range("D95").Value = 41
$ws.Range("E95").Value = 439
$ws.Range("H95").Value = 11

# --- Rows 122/123: Martinica / Guatemala swap ranking ---
# Guatemala overtakes Martinica: new row 122 = Guatemala (fresh numbers),
# new row 123 = Martinica (carries forward the old row-122 numbers unchanged).
$ws.Range("A122").Value = "Guatemala"
$ws.Range("C122").Value = 18
$ws.Range("D122").Value = 19
$ws.Range("E122").Value = 131
$ws.Range("F122").Value = 3
$ws.Range("G122").Value = 2
$ws.Range("H122").Value = 5

$ws.Range("A123").Value = "Martinica"
$ws.Range("B123").Value = 155
$ws.Range("C123").Value = 0
$ws.Range("D123").Value = 50
$ws.Range("E123").Value = 99
$ws.Range("F123").Value = 19
$ws.Range("H123").Value = 6

# --- Row 179 (Dominica) ---
$ws.Range("D179").Value = 7
$ws.Range("E179").Value = 9
